$d = $word.ActiveDocument
$d.Content.Find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
